$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content corrections (contacts sheet) ---
# Row 2 contact: first name typo fix, email typo fix, trailing separator
# removed from hobbies list.
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 contact: first name + street updated, trailing separator removed
# from hobbies list.
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Formatting touch-up ---
# Slightly taller header/data rows.
$ws.Rows("1:3").RowHeight = 19.5
